$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.100.13"
$ws.Range("E2").Value = "  -2.08%  "

$ws.Range("D3").Value = "1.802.92"
$ws.Range("E3").Value = "  -2.27%  "

$cell = $ws.Range("D4")
$cell.Value = "'1.006"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.39%  "

$cell = $ws.Range("D5")
$cell.Value = "'309.11"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -2.04%  "

$ws.Range("E6").Value = "  +0.24%  "

$cell = $ws.Range("D7")
$cell.Value = "'0.4236"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -2.01%  "

$cell = $ws.Range("D8")
$cell.Value = "'0.3609"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -2.66%  "

$cell = $ws.Range("D9")
$cell.Value = "'0.07240"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -1.37%  "

$cell = $ws.Range("D10")
$cell.Value = "'0.8454"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -3.78%  "

$cell = $ws.Range("D11")
$cell.Value = "'20.36"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -3.30%  "

$ws.Range("D12").Value = "1.788.24"
$ws.Range("E12").Value = "  -4.59%  "

$cell = $ws.Range("D13")
$cell.Value = "'5.298"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -3.26%  "

$cell = $ws.Range("D14")
$cell.Value = "'6.380"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -3.34%  "

$cell = $ws.Range("D15")
$cell.Value = "'0.06795"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -2.30%  "

$cell = $ws.Range("D16")
$cell.Value = "'1.003"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -0.01%  "

$cell = $ws.Range("D17")
$cell.Value = "'81.00"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -0.21%  "

$cell = $ws.Range("D18")
$cell.Value = "'0.000008765"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -3.31%  "

$cell = $ws.Range("D19")
$cell.Value = "'1.005"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +0.30%  "

$cell = $ws.Range("D20")
$cell.Value = "'15.07"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -3.38%  "

$ws.Range("D21").Value = "27.326.97"
$ws.Range("E21").Value = "  -1.80%  "

$cell = $ws.Range("D22")
$cell.Value = "'5.092"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -0.72%  "

$cell = $ws.Range("D23")
$cell.Value = "'11.11"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.88%  "

$ws.Range("D24").Value = "2.084.17"
$ws.Range("E24").Value = "  -3.21%  "

$ws.Range("E25").Value = "  -1.70%  "

$cell = $ws.Range("D26")
$cell.Value = "'153.14"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -0.61%  "

$cell = $ws.Range("D27")
$cell.Value = "'18.23"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -3.72%  "

$cell = $ws.Range("D28")
$cell.Value = "'5.035"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -5.26%  "

$cell = $ws.Range("D29")
$cell.Value = "'114.13"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -1.61%  "

$cell = $ws.Range("D30")
$cell.Value = "'1.660"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -11.71%  "

$cell = $ws.Range("D31")
$cell.Value = "'0.09008"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +0.83%  "

$cell = $ws.Range("D32")
$cell.Value = "'0.7368"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -6.56%  "

$cell = $ws.Range("D33")
$cell.Value = "'2.884"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -3.12%  "

$cell = $ws.Range("D34")
$cell.Value = "'4.370"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -5.41%  "

$cell = $ws.Range("D35")
$cell.Value = "'1.099"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -6.88%  "

$ws.Range("E36").Value = "  +0.20%  "

$cell = $ws.Range("D37")
$cell.Value = "'1.082"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -2.08%  "

$cell = $ws.Range("D38")
$cell.Value = "'0.05157"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -5.19%  "

$cell = $ws.Range("D39")
$cell.Value = "'0.01909"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -2.76%  "

$cell = $ws.Range("D40")
$cell.Value = "'0.1636"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -3.21%  "

$cell = $ws.Range("D41")
$cell.Value = "'0.4991"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -3.55%  "

$cell = $ws.Range("D42")
$cell.Value = "'2.619"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -8.10%  "

$cell = $ws.Range("D43")
$cell.Value = "'8.140"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -5.93%  "

$cell = $ws.Range("D44")
$cell.Value = "'5.963"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -12.41%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$cell = $ws.Range("D45")
$cell.Value = "'10.31"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -3.57%  "

$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$cell = $ws.Range("D46")
$cell.Value = "'105.15"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -1.47%  "

$ws.Range("E47").Value = "  +0.23%  "

$cell = $ws.Range("D48")
$cell.Value = "'0.06328"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -3.48%  "

$cell = $ws.Range("D49")
$cell.Value = "'0.4546"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -4.93%  "

$cell = $ws.Range("D50")
$cell.Value = "'1.606"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -3.56%  "

$ws.Range("E51").Value = "  -6.47%  "
